$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Serial 1 - Frontslash
$ws.Cells.Item(2, 2).Value = 'Frontslash'
$ws.Cells.Item(2, 3).Value = 'Ankush Gautam, Vaibhav Srivastva, utkal, Aishlee Joshi'
$ws.Cells.Item(2, 9).Value = 'MLSC274581924053'
$ws.Cells.Item(2, 10).Value = 2985

# Row 3: Serial 2 - xantiedar
$ws.Cells.Item(3, 2).Value = 'xantiedar'
$ws.Cells.Item(3, 3).Value = 'Deb, Tanish Gupta, Achin, Nidhi'
$ws.Cells.Item(3, 9).Value = 'MLSC273411206789'
$ws.Cells.Item(3, 10).Value = 1618

# Row 4: Serial 3 - BCS
$ws.Cells.Item(4, 2).Value = 'BCS'
$ws.Cells.Item(4, 3).Value = 'Siddharth, Geet, Yashit Arora, Aryan Thakkar'
$ws.Cells.Item(4, 9).Value = 'MLSC278956012348'
$ws.Cells.Item(4, 10).Value = 3875

# Row 5: Serial 4 - Home_Team
$ws.Cells.Item(5, 2).Value = 'Home_Team'
$ws.Cells.Item(5, 3).Value = 'Akanksha, Hurreet, Jasman, Puranjay'
$ws.Cells.Item(5, 9).Value = 'MLSC271900439281'
$ws.Cells.Item(5, 10).Value = 2930

# Row 6: Serial 5 - GajarKaHalwa
$ws.Cells.Item(6, 2).Value = 'GajarKaHalwa'
$ws.Cells.Item(6, 3).Value = 'Harshit, Nikunj Dewan, Tulika, Pranav'
$ws.Cells.Item(6, 9).Value = 'MLSC276753908823'
$ws.Cells.Item(6, 10).Value = 4347

# Row 7: Serial 6 - x-tasy
$ws.Cells.Item(7, 2).Value = 'x-tasy'
$ws.Cells.Item(7, 3).Value = 'Sameer Khan, Karan, Vinesh, Shorya'
$ws.Cells.Item(7, 9).Value = 'MLSC278021677349'
$ws.Cells.Item(7, 10).Value = 7282

# Row 8: Serial 7 - BlackHole
$ws.Cells.Item(8, 2).Value = 'BlackHole'
$ws.Cells.Item(8, 3).Value = 'Tanay, Rohin, Rudra, Krish'
$ws.Cells.Item(8, 9).Value = 'MLSC279188325690'
$ws.Cells.Item(8, 10).Value = 5103

# Row 9: Serial 8 - Vanguard
$ws.Cells.Item(9, 2).Value = 'Vanguard'
$ws.Cells.Item(9, 3).Value = 'Suryansh, Vaibhav, mankirat, Saket'
$ws.Cells.Item(9, 9).Value = 'MLSC274012093948'
$ws.Cells.Item(9, 10).Value = 5387

# Row 10: Serial 9 - Andhadun Players
$ws.Cells.Item(10, 2).Value = 'Andhadun Players'
$ws.Cells.Item(10, 3).Value = 'Amanjot, Jahanvi, Ram, Dhiren'
$ws.Cells.Item(10, 9).Value = 'MLSC273665718204'
$ws.Cells.Item(10, 10).Value = 5935

# Row 11: Serial 10 - BlueBull
$ws.Cells.Item(11, 2).Value = 'BlueBull'
$ws.Cells.Item(11, 3).Value = 'Shreyansh, Parvesh Lamba, Shriyam, Tanveer'
$ws.Cells.Item(11, 9).Value = 'MLSC279937456132'
$ws.Cells.Item(11, 10).Value = 2268

# Row 12: Serial 11 - GenF
$ws.Cells.Item(12, 2).Value = 'GenF'
$ws.Cells.Item(12, 3).Value = 'Piyush, Chirag, Bhagya, Vandini Garg'
$ws.Cells.Item(12, 9).Value = 'MLSC275302947685'
$ws.Cells.Item(12, 10).Value = 3591

# Row 13: Serial 12 - Jaguar
$ws.Cells.Item(13, 2).Value = 'Jaguar'
$ws.Cells.Item(13, 3).Value = 'Anushka Verma, Manvi, Daksh Aggarwal, Uday'
$ws.Cells.Item(13, 9).Value = 'MLSC272490411236'
$ws.Cells.Item(13, 10).Value = 5954

# Row 14: Serial 13 - Guns and Roses
$ws.Cells.Item(14, 2).Value = 'Guns and Roses'
$ws.Cells.Item(14, 3).Value = 'Prithvi, Sahil, Vrishank, Vani'
$ws.Cells.Item(14, 9).Value = 'MLSC278386074821'
$ws.Cells.Item(14, 10).Value = 2457

# Row 15: Serial 14 - Homies
$ws.Cells.Item(15, 2).Value = 'Homies'
$ws.Cells.Item(15, 3).Value = 'Siddhant, Aryan Thakur, Chaitanya, Prerna Garg'
$ws.Cells.Item(15, 9).Value = 'MLSC271219486573'
$ws.Cells.Item(15, 10).Value = 2740

# Row 16: Serial 15 - BABLU
$ws.Cells.Item(16, 2).Value = 'BABLU'
$ws.Cells.Item(16, 3).Value = 'Vansh, Harsh, Shikhar, Vaibhav Sundriyal'
$ws.Cells.Item(16, 8).Value = '5, 2, 7'
$ws.Cells.Item(16, 9).Value = 'MLSC275630089147'
$ws.Cells.Item(16, 10).Value = 401

# Row 17: Serial 16 - OUTLAW
$ws.Cells.Item(17, 2).Value = 'OUTLAW'
$ws.Cells.Item(17, 3).Value = 'Arshit, ayush, pragun, aarav'
$ws.Cells.Item(17, 9).Value = 'MLSC277953712340'
$ws.Cells.Item(17, 10).Value = 1985

# Row 18: Serial 17 - SQUARE
$ws.Cells.Item(18, 2).Value = 'SQUARE'
$ws.Cells.Item(18, 3).Value = 'Jasmine, Riya, jyotsna, Anushka'
$ws.Cells.Item(18, 9).Value = 'MLSC273519849023'
$ws.Cells.Item(18, 10).Value = 0

# Row 19: Serial 18 - TSM_entity
$ws.Cells.Item(19, 2).Value = 'TSM_entity'
$ws.Cells.Item(19, 3).Value = 'shounok, amit, aryan, aditya sharma'
$ws.Cells.Item(19, 9).Value = 'MLSC272764021980'
$ws.Cells.Item(19, 10).Value = 1985

# Row 20: Serial 19 - Maqsad Returns
$ws.Cells.Item(20, 2).Value = 'Maqsad Returns'
$ws.Cells.Item(20, 3).Value = 'Yash, Dixant, Ankit Pandey, vrattan'
$ws.Cells.Item(20, 9).Value = 'MLSC278241857304'
$ws.Cells.Item(20, 10).Value = 4309

# Row 21: Serial 20 - Syrups
$ws.Cells.Item(21, 2).Value = 'Syrups'
$ws.Cells.Item(21, 3).Value = 'Aashish, Abhishek, Ujjwal, Vishesh'
$ws.Cells.Item(21, 9).Value = 'MLSC279401358492'
$ws.Cells.Item(21, 10).Value = 3875

# Row 22: Serial 21 - Team Holmes
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 'Team Holmes'
$ws.Cells.Item(22, 3).Value = 'Aayushman, Madhav Gaba, Saksham Katna, hemant'
$ws.Cells.Item(22, 4).Value = '3, 1, 2'
$ws.Cells.Item(22, 5).Value = '101, 102, 103, 104, 105'
$ws.Cells.Item(22, 6).Value = '201, 202, 203, 204, 205'
$ws.Cells.Item(22, 7).Value = '301, 302, 303, 304, 305'
$ws.Cells.Item(22, 8).Value = "'2"
$ws.Cells.Item(22, 9).Value = 'MLSC275146789013'
$ws.Cells.Item(22, 10).Value = 369

# Row 23: Serial 22 - CRESTFALL
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = 'CRESTFALL'
$ws.Cells.Item(23, 3).Value = 'Aditya Gupta, Raghav, Aarush, Amish'
$ws.Cells.Item(23, 4).Value = '1, 2, 3'
$ws.Cells.Item(23, 5).Value = '102, 103, 104, 105, 101'
$ws.Cells.Item(23, 6).Value = '202, 203, 204, 205, 201'
$ws.Cells.Item(23, 7).Value = '302, 303, 304, 305, 301'
$ws.Cells.Item(23, 8).Value = ''
$ws.Cells.Item(23, 9).Value = 'MLSC273805276149'
$ws.Cells.Item(23, 10).Value = 6520

# Row 24: Serial 23 - Chdi Gang
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = 'Chdi Gang'
$ws.Cells.Item(24, 3).Value = 'Sameer Verma, Aditya, Palak, Bhavninder'
$ws.Cells.Item(24, 4).Value = '2, 3, 1'
$ws.Cells.Item(24, 5).Value = '103, 104, 105, 101, 102'
$ws.Cells.Item(24, 6).Value = '203, 204, 205, 201, 202'
$ws.Cells.Item(24, 7).Value = '303, 304, 305, 301, 302'
$ws.Cells.Item(24, 8).Value = ''
$ws.Cells.Item(24, 9).Value = 'MLSC271729503826'
$ws.Cells.Item(24, 10).Value = 3780

# Row 25: Serial 24 - Om Rajpal
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = 'Om Rajpal'
$ws.Cells.Item(25, 3).Value = 'Aman, Angad, harshil, Kashish'
$ws.Cells.Item(25, 4).Value = '3, 1, 2'
$ws.Cells.Item(25, 5).Value = '104, 105, 101, 102, 103'
$ws.Cells.Item(25, 6).Value = '204, 205, 201, 202, 203'
$ws.Cells.Item(25, 7).Value = '304, 305, 301, 302, 303'
$ws.Cells.Item(25, 8).Value = ''
$ws.Cells.Item(25, 9).Value = 'MLSC276089314578'
$ws.Cells.Item(25, 10).Value = 3100

# Row 26: Serial 25 - HR10
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = 'HR10'
$ws.Cells.Item(26, 3).Value = 'Sourav, Ishtpreet, Yashkaran, Agrim Bhatt'
$ws.Cells.Item(26, 4).Value = '1, 2, 3'
$ws.Cells.Item(26, 5).Value = '105, 101, 102, 103, 104'
$ws.Cells.Item(26, 6).Value = '205, 201, 202, 203, 204'
$ws.Cells.Item(26, 7).Value = '305, 301, 302, 303, 304'
$ws.Cells.Item(26, 8).Value = ''
$ws.Cells.Item(26, 9).Value = 'MLSC274920348612'
$ws.Cells.Item(26, 10).Value = 662

# Row 27: Serial 26 - Thalaforareason
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = 'Thalaforareason'
$ws.Cells.Item(27, 3).Value = 'Aakarsh, Aryan2, Anirudh, Pranav7'
$ws.Cells.Item(27, 4).Value = '2, 3, 1'
$ws.Cells.Item(27, 5).Value = '101, 102, 103, 104, 105'
$ws.Cells.Item(27, 6).Value = '201, 202, 203, 204, 205'
$ws.Cells.Item(27, 7).Value = '301, 302, 303, 304, 305'
$ws.Cells.Item(27, 8).Value = ''
$ws.Cells.Item(27, 9).Value = 'MLSC273276041398'
$ws.Cells.Item(27, 10).Value = 2169

# Row 28: Serial 27 - Reapers
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 'Reapers'
$ws.Cells.Item(28, 3).Value = 'Ishwinder, Aradhya, Kumud, Arnav'
$ws.Cells.Item(28, 4).Value = '3, 1, 2'
$ws.Cells.Item(28, 5).Value = '102, 103, 104, 105, 101'
$ws.Cells.Item(28, 6).Value = '202, 203, 204, 205, 201'
$ws.Cells.Item(28, 7).Value = '302, 303, 304, 305, 301'
$ws.Cells.Item(28, 8).Value = ''
$ws.Cells.Item(28, 9).Value = 'MLSC275490028347'
$ws.Cells.Item(28, 10).Value = 2741
